$p = $ppt.ActivePresentation

# --- Slide 22, shape "CuadroTexto 2" ---------------------------------------
# "en cada colectivo y nivel de dirección." -> split into two runs:
#   "en cada colectivo y nivel de " + "dirección."
$s22 = $p.Slides.Item(22)
$shp22a = $s22.Shapes.Item(2)
$tr22a = $shp22a.TextFrame.TextRange
$para = $tr22a.Paragraphs(2)
$run = $para.Runs(3)
$run.Text = "en cada colectivo y nivel de "
[void]$run.InsertAfter("dirección.")

# --- Slide 22, shape "CuadroTexto 4" ---------------------------------------
# "...para cada curso académico. " ->
#   "...para cada curso " + "académico" + "."
# (the long run also loses its trailing space in the process)
$shp22b = $s22.Shapes.Item(4)
$tr22b = $shp22b.TextFrame.TextRange
$fullText = $tr22b.Text
$marker = "para cada curso académico. "
$idx = $fullText.IndexOf($marker)
$prefixLen = "para cada curso ".Length
$acadStart = $idx + $prefixLen + 1
$acadRange = $tr22b.Characters($acadStart, "académico".Length)
$acadRange.Text = "académico"
$dotStart = $acadStart + "académico".Length
$dotRange = $tr22b.Characters($dotStart, 2)
$dotRange.Text = "."

# --- Slide 4, shape "CuadroTexto 2" -----------------------------------------
# Merge the four runs "Se " + "llevará a cabo en cada uno de " +
# "los niveles " + "organizativos del proceso docente educativo."
# into a single run.
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange
$fullText4 = $tr4.Text
$idx4 = $fullText4.IndexOf("Se llevará a cabo")
$target4 = $tr4.Characters($idx4 + 1, $fullText4.Length - $idx4)
$target4.Text = "Se llevará a cabo en cada uno de los niveles organizativos del proceso docente educativo."
